$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: (301, 6, 45, 30, 60, 45) -> (801, 3, 67, 65, 52, 45)
$ws.Range("A2").Value = 801
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 67
$ws.Range("D2").Value = 65
$ws.Range("E2").Value = 52

# Row 3: (701, 3, 90, 45, 97, 15) -> (1203, 3, 15, 15, 15, 15)
$ws.Range("A3").Value = 1203
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 15

# Row 4: (1203, 3, 15, 15, 15, 15) -> (101, 9, 30, 15, 60, 15)
$ws.Range("A4").Value = 101
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 30
$ws.Range("E4").Value = 60

# Row 5: (101, 9, 30, 15, 60, 15) -> (401, 9, 48, 67, 75, 45)
$ws.Range("A5").Value = 401
$ws.Range("C5").Value = 48
$ws.Range("D5").Value = 67
$ws.Range("E5").Value = 75
$ws.Range("F5").Value = 45

# Row 6: (901, 16, 15, 45, 60, 60) -> (1201, 2, 10, 10, 10, 10)
$ws.Range("A6").Value = 1201
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 10

# Row 7: (201, 9, 30, 15, 45, 30) -> (501, 9, 52, 30, 75, 45)
$ws.Range("A7").Value = 501
$ws.Range("C7").Value = 52
$ws.Range("D7").Value = 30
$ws.Range("E7").Value = 75
$ws.Range("F7").Value = 45

# Row 8: (1001, 18, 30, 75, 60, 72) -> (701, 3, 90, 45, 97, 15)
$ws.Range("A8").Value = 701
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 90
$ws.Range("D8").Value = 45
$ws.Range("E8").Value = 97
$ws.Range("F8").Value = 15

# Row 9: (501, 9, 52, 30, 75, 45) -> (601, 9, 60, 67, 60, 42)
$ws.Range("A9").Value = 601
$ws.Range("C9").Value = 60
$ws.Range("D9").Value = 67
$ws.Range("E9").Value = 60
$ws.Range("F9").Value = 42

# Row 10: (1202, 2, 10, 10, 10, 10) -> (902, 1, 0, 0, 0, 0)
$ws.Range("A10").Value = 902
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0

# Row 11: (601, 9, 60, 67, 60, 42) -> (1001, 18, 30, 75, 60, 72)
$ws.Range("A11").Value = 1001
$ws.Range("B11").Value = 18
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 75
$ws.Range("F11").Value = 72

# Row 12: (801, 3, 67, 65, 52, 45) -> (301, 6, 45, 30, 60, 45)
$ws.Range("A12").Value = 301
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 45
$ws.Range("D12").Value = 30
$ws.Range("E12").Value = 60

# Row 13: (1201, 2, 10, 10, 10, 10) -> (901, 16, 15, 45, 60, 60)
$ws.Range("A13").Value = 901
$ws.Range("B13").Value = 16
$ws.Range("C13").Value = 15
$ws.Range("D13").Value = 45
$ws.Range("E13").Value = 60
$ws.Range("F13").Value = 60

# Row 14: (902, 1, 0, 0, 0, 0) -> (201, 9, 30, 15, 45, 30)
$ws.Range("A14").Value = 201
$ws.Range("B14").Value = 9
$ws.Range("C14").Value = 30
$ws.Range("D14").Value = 15
$ws.Range("E14").Value = 45
$ws.Range("F14").Value = 30

# Row 15: (401, 9, 48, 67, 75, 45) -> (1202, 2, 10, 10, 10, 10)
$ws.Range("A15").Value = 1202
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 10

# Row 16: (802, 0, 4, 5, 4, 0) -> (3, 0, 3, 3, 3, 3)
$ws.Range("A16").Value = 3
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 3

# Row 17: (1, 0, 2, 2, 2, 2) -> (2, 0, 2, 2, 2, 2)
$ws.Range("A17").Value = 2

# Row 18: (3, 0, 3, 3, 3, 3) -> (502, 0, 4, 0, 0, 0)
$ws.Range("A18").Value = 502
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0

# Row 19: (502, 0, 4, 0, 0, 0) -> (1101, 0, 15, 30, 30, 0)
$ws.Range("A19").Value = 1101
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 30
$ws.Range("E19").Value = 30

# Row 20: (1101, 0, 15, 30, 30, 0) -> (1, 0, 2, 2, 2, 2)
$ws.Range("A20").Value = 1
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 2

# Row 21: (2, 0, 2, 2, 2, 2) -> (802, 0, 4, 5, 4, 0)
$ws.Range("A21").Value = 802
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 0

# Row 22: (602, 0, 0, 4, 0, 9) -> (402, 0, 0, 4, 0, 0)
$ws.Range("A22").Value = 402
$ws.Range("F22").Value = 0

# Row 23: (402, 0, 0, 4, 0, 0) -> (602, 0, 0, 4, 0, 9)
$ws.Range("A23").Value = 602
$ws.Range("F23").Value = 9

Write-Output "Applied row re-ordering changes"
